$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell above the table: "Энергопотребление" -> "Энергопотребление (мА)"
$ws.Range("A9").Value = "Энергопотребление (мА)"

# --- Table header row 10 (unchanged text, just keeping values explicit / no-op safe)
$ws.Range("B10").Value = "Сон"
$ws.Range("C10").Value = "Рабочее"
$ws.Range("D10").Value = "Время работы"

# --- Row 11: STM32F103 - real measured resistor-divider currents (numbers, not text)
$ws.Range("A11").Value = "STM32F103"
$ws.Range("B11").Value = 5.2
$ws.Range("C11").Value = 8.4
$ws.Range("D11").Formula = "=D15"

# --- Row 12: BME280
$ws.Range("A12").Value = "BME280"
$ws.Range("B12").Value = 0.07
$ws.Range("C12").Value = 0.2
$ws.Range("D12").Formula = "=D15"

# --- Row 13: ESP8266
$ws.Range("A13").Value = "ESP8266"
$ws.Range("B13").Value = 0.3
$ws.Range("C13").Value = 70
$ws.Range("D13").Formula = "=D15"
# D13 picks up D15's number-format style (border/fill) instead of the row's own style
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122)

# --- Row 14: Стабилизатор
$ws.Range("A14").Value = "Стабилизатор"
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = "всегда"

# --- Row 15: В сборе (totals)
$ws.Range("A15").Value = "В сборе"
$ws.Range("B15").Formula = "=B11+B12+B13+B14"
$ws.Range("C15").Formula = "=C11+C12+C13+C14"
$ws.Range("D15").Value = "6 с"

# --- Row 17: expected runtime without recharge, updated numbers
$ws.Range("A17").Value = "Ожидаемое время работы без подзарядки при батарее 2А/ч = 8,3 дней"

# --- Row 18 (new row of content): average consumption note, formatted/merged like row 17
$ws.Range("A18:D18").Merge()
$ws.Range("A18").Value = "Среднее энергопотребление = 10 мА"
$ws.Range("A18:D18").HorizontalAlignment = -4108
$ws.Range("A18:D18").VerticalAlignment = -4108

# --- Update selection / active cell to match the saved view state
$ws.Range("D19").Select()
